$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.379.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "'1.830.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'317.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "'0.5347"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'0.4063"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.68%  "
$ws.Range("D9").Value = "'0.07611"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("D10").Value = "'41.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").Value = "'1.106"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").Value = "'6.355"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.43%  "
$ws.Range("D13").Value = "'1.002"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "'7.588"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.15%  "
$ws.Range("D15").Value = "'20.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").Value = "'1.827.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.02%  "
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "'89.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").Value = "'0.06624"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("D20").Value = "'17.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("D23").Value = "'28.400.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("D24").Value = "'11.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").Value = "'2.171"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.93%  "
$ws.Range("D26").Value = "'2.477"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.41%  "
$ws.Range("D27").Value = "'157.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").Value = "'2.040.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.22%  "
$ws.Range("D30").Value = "'124.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.53%  "
$ws.Range("D31").Value = "'1.121"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.47%  "
$ws.Range("D32").Value = "'0.1094"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.84%  "
$ws.Range("D33").Value = "'5.688"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.43%  "
$ws.Range("D34").Value = "'3.628"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("D35").Value = "'0.07164"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.67%  "
$ws.Range("D36").Value = "'0.2260"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").Value = "'0.02343"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.59%  "
$ws.Range("D38").Value = "'5.216"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.92%  "
$ws.Range("D39").Value = "'8.835"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.53%  "
$ws.Range("D40").Value = "'0.6276"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.92%  "
$ws.Range("D41").Value = "'11.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.86%  "
$ws.Range("D42").Value = "'1.190"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "'1.400"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("D45").Value = "'13.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.45%  "
$ws.Range("D46").Value = "'3.707"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("D47").Value = "'0.5856"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").Value = "'125.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "'1.993"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.49%  "
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").Value = "'0.06897"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.08%  "
